$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.474.07'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.839.89'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''319.17'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = '''0.5323'
$ws.Range('E7').Value = '  -1.55%  '
$ws.Range('D8').Value = '''0.4027'
$ws.Range('E8').Value = '  +6.51%  '
$ws.Range('D9').Value = '''0.07608'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').Value = '''41.84'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '''1.109'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '''6.332'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''7.634'
$ws.Range('E13').Value = '  +4.63%  '
$ws.Range('B14').Value = 'BinanceUSD'
$ws.Range('C14').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D14').Value = '''0.9989'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '''20.82'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '1.833.04'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('D17').Value = '''90.00'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').Value = '''0.00001075'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''0.06610'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('D20').Value = '''17.74'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '''0.9991'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = '''6.081'
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('D23').Value = '28.465.23'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''11.22'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').Value = '''2.109'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''20.68'
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''156.94'
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''2.442'
$ws.Range('E28').Value = '  +5.10%  '
$ws.Range('D29').Value = '2.042.45'
$ws.Range('E29').Value = '  +2.08%  '
$ws.Range('D30').Value = '''123.90'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = '''1.117'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('D32').Value = '''0.1099'
$ws.Range('E32').Value = '  +4.30%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.690'
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.663'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').Value = '''0.07250'
$ws.Range('E35').Value = '  +12.26%  '
$ws.Range('D36').Value = '''0.2256'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '''5.267'
$ws.Range('E37').Value = '  +4.88%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02351'
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('D39').Value = '''8.818'
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('D40').Value = '''11.38'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').Value = '''0.6294'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '''1.203'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('D44').Value = '''0.9990'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('D46').Value = '''3.707'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('D47').Value = '''0.5849'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').Value = '''126.01'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').Value = '''1.983'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('D50').Value = '''1.197'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').Value = '''0.06922'
$ws.Range('E51').Value = '  +0.52%  '
